$wb = $excel.ActiveWorkbook

# --- Update the conversion message on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.14 = 28808.17 pesos`n✅ 28808.17 pesos = 7.12 = 965.73 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the rate values on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 139.995
$wsTasas.Range("O10").Value = 4033
$wsTasas.Range("N12").Value = 4045
$wsTasas.Range("O12").Value = 135.6
